$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates --------------------------------------------
# These values are authored as literal text (e.g. "30.590.53", "1.001",
# "246.00") rather than numbers, so each cell is force-formatted as Text
# before the value is written - otherwise Excel's normal text-to-number
# coercion would turn things like "1.001" or "246.00" into the numbers
# 1.001 / 246. The style is reset back to "Normal" afterwards so no stray
# cell formatting is left behind.
$priceUpdates = @{
    "D2"  = "30.590.53"
    "D3"  = "1.919.62"
    "D5"  = "246.00"
    "D6"  = "1.001"
    "D8"  = "0.2913"
    "D9"  = "0.06724"
    "D10" = "111.58"
    "D11" = "19.39"
    "D12" = "1.924.66"
    "D13" = "0.07589"
    "D14" = "5.342"
    "D15" = "0.6741"
    "D16" = "294.68"
    "D17" = "30.607.94"
    "D18" = "13.09"
    "D19" = "1.001"
    "D20" = "0.000007569"
    "D21" = "2.177.10"
    "D22" = "5.529"
    "D24" = "6.428"
    "D25" = "9.487"
    "D26" = "164.98"
    "D27" = "20.32"
    "D28" = "2.108"
    "D30" = "1.438"
    "D31" = "4.138"
    "D32" = "4.102"
    "D33" = "0.05030"
    "D34" = "0.7421"
    "D35" = "1.141"
    "D36" = "0.9999"
    "D37" = "0.02031"
    "D38" = "2.705"
    "D39" = "2.689"
    "D40" = "2.025"
    "D41" = "109.99"
    "D42" = "0.4460"
    "D43" = "0.8666"
    "D44" = "5.879"
    "D45" = "70.05"
    "D46" = "1.001"
    "D47" = "7.273"
    "D48" = "48.43"
    "D49" = "9.276"
    "D50" = "0.1231"
    "D51" = "0.2543"
}

foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}

# --- Column E ("Volume(1h)") updates ----------------------------------------
$volumeUpdates = @{
    "E2"  = "  +0.09%  "
    "E3"  = "  +0.03%  "
    "E4"  = "  -0.01%  "
    "E5"  = "  -0.50%  "
    "E6"  = "  +0.01%  "
    "E7"  = "  +2.72%  "
    "E8"  = "  +1.28%  "
    "E9"  = "  -1.21%  "
    "E10" = "  +6.32%  "
    "E11" = "  +5.57%  "
    "E12" = "  +0.27%  "
    "E13" = "  -1.38%  "
    "E15" = "  +0.95%  "
    "E16" = "  +2.26%  "
    "E17" = "  +0.11%  "
    "E18" = "  +1.29%  "
    "E19" = "  +0.05%  "
    "E20" = "  -0.27%  "
    "E21" = "  +0.53%  "
    "E22" = "  -0.26%  "
    "E23" = "  -0.09%  "
    "E24" = "  +2.11%  "
    "E25" = "  +1.20%  "
    "E26" = "  -2.17%  "
    "E27" = "  -3.90%  "
    "E28" = "  -0.47%  "
    "E29" = "  +0.67%  "
    "E30" = "  +3.08%  "
    "E31" = "  -0.82%  "
    "E32" = "  +0.55%  "
    "E33" = "  -0.10%  "
    "E34" = "  +0.69%  "
    "E35" = "  -0.44%  "
    "E37" = "  -2.11%  "
    "E38" = "  -1.61%  "
    "E39" = "  +0.11%  "
    "E40" = "  -1.86%  "
    "E41" = "  -1.20%  "
    "E42" = "  +1.50%  "
    "E43" = "  -1.47%  "
    "E44" = "  +0.06%  "
    "E45" = "  +4.25%  "
    "E46" = "  +0.05%  "
    "E47" = "  +0.21%  "
    "E48" = "  +0.31%  "
    "E49" = "  -0.58%  "
    "E50" = "  +0.08%  "
    "E51" = "  +3.26%  "
}

foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}

# --- Rows 40/41 swapped places (RenderToken now ranks above Quant) ---------
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
